# Update the Sem-1 timetable: Section_A and Section_B sheets get new
# course-code assignments for the 9:00-10:30 / 10:30-12:00 / 14:00-15:30 /
# 15:30-17:00 / 17:00-18:30 slots (Mon-Fri). The 12:30-14:00 "LUNCH BREAK"
# row and the header row/column are left untouched.

$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("Section_A")
$wsB = $wb.Worksheets.Item("Section_B")

# ---- Section_A ----
$wsA.Range("B2").Value = "MA102"
$wsA.Range("C2").Value = "CS101"
$wsA.Range("D2").Value = "EC101"
$wsA.Range("E2").Value = "Free"
$wsA.Range("F2").Value = "MA101"

$wsA.Range("B3").Value = "CS151"
$wsA.Range("C3").Value = "CS151"
$wsA.Range("D3").Value = "Free"
$wsA.Range("E3").Value = "HS101"
$wsA.Range("F3").Value = "Free"

$wsA.Range("B5").Value = "HS101"
$wsA.Range("C5").Value = "Free"
$wsA.Range("D5").Value = "Free"
$wsA.Range("E5").Value = "DS101"
$wsA.Range("F5").Value = "MA102"

$wsA.Range("B6").Value = "CS101"
$wsA.Range("C6").Value = "Free"
$wsA.Range("D6").Value = "Free"
$wsA.Range("E6").Value = "CS101"
$wsA.Range("F6").Value = "DS101"

$wsA.Range("B7").Value = "Free"
$wsA.Range("C7").Value = "EC101"
$wsA.Range("D7").Value = "HS101"
$wsA.Range("E7").Value = "MA101"
$wsA.Range("F7").Value = "EC101"

# ---- Section_B ----
$wsB.Range("B2").Value = "MA102"
$wsB.Range("C2").Value = "Free"
$wsB.Range("D2").Value = "Free"
$wsB.Range("E2").Value = "EC101"
$wsB.Range("F2").Value = "MA101"

$wsB.Range("B3").Value = "MA101"
$wsB.Range("C3").Value = "HS101"
$wsB.Range("D3").Value = "CS101"
$wsB.Range("E3").Value = "MA102"
$wsB.Range("F3").Value = "CS101"

$wsB.Range("B5").Value = "CS151"
$wsB.Range("C5").Value = "Free"
$wsB.Range("D5").Value = "Free"
$wsB.Range("E5").Value = "HS101"
$wsB.Range("F5").Value = "DS101"

$wsB.Range("B6").Value = "Free"
$wsB.Range("C6").Value = "EC101"
$wsB.Range("D6").Value = "HS101"
$wsB.Range("E6").Value = "DS101"
$wsB.Range("F6").Value = "CS151"

$wsB.Range("B7").Value = "EC101"
$wsB.Range("C7").Value = "CS101"
$wsB.Range("D7").Value = "Free"
$wsB.Range("E7").Value = "Free"
$wsB.Range("F7").Value = "Free"
